# Append an extra sentence to the "கர்ம தோஷத்தை..." explanation in cell B10.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("B10")
$addition = "இல்ல தெய்வம் அனுக்கிரகம் பற்றி கூறும் இடம்,பொருளாதார உயர்வு,அதிகாரத்தை சொல்லும் இடம்."
$cell.Value = $cell.Value() + $addition

# The longer text now needs a taller row, while row 3 settles back down.
$ws.Rows.Item(3).RowHeight = 60
$ws.Rows.Item(10).RowHeight = 75

# Editing B10 and pressing Enter leaves the selection on the next cell down.
$ws.Range("B11").Select()
